$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) and "Volume(1h)" (column E) figures for the
# symbol list refresh, matching the author's source-data re-scrape.
$updates = @(
    @{ Cell = "D2"; Value = "323.40" },
    @{ Cell = "E2"; Value = "8.78%" },
    @{ Cell = "D3"; Value = "49.62" },
    @{ Cell = "E3"; Value = "18.73%" },
    @{ Cell = "D4"; Value = "5.294" },
    @{ Cell = "E4"; Value = "5.48%" },
    @{ Cell = "D5"; Value = "0.08161" },
    @{ Cell = "E5"; Value = "8.37%" },
    @{ Cell = "D6"; Value = "4.595" },
    @{ Cell = "E6"; Value = "5.02%" },
    @{ Cell = "D7"; Value = "1.684" },
    @{ Cell = "E7"; Value = "3.94%" },
    @{ Cell = "D8"; Value = "1.213" },
    @{ Cell = "E8"; Value = "31.62%" },
    @{ Cell = "D9"; Value = "0.1351" },
    @{ Cell = "E9"; Value = "14.29%" },
    @{ Cell = "D10"; Value = "0.1968" },
    @{ Cell = "E10"; Value = "7.47%" },
    @{ Cell = "D11"; Value = "0.09702" },
    @{ Cell = "E11"; Value = "7.21%" },
    @{ Cell = "D12"; Value = "0.04418" },
    @{ Cell = "E12"; Value = "8.12%" },
    @{ Cell = "E13"; Value = "-0.24%" },
    @{ Cell = "D14"; Value = "0.001327" },
    @{ Cell = "E14"; Value = "3.78%" },
    @{ Cell = "D15"; Value = "0.005804" },
    @{ Cell = "E15"; Value = "-0.31%" },
    @{ Cell = "E16"; Value = "1.07%" },
    @{ Cell = "D17"; Value = "2.438" },
    @{ Cell = "E17"; Value = "1.55%" },
    @{ Cell = "E18"; Value = "2.01%" },
    @{ Cell = "D19"; Value = "8.160" },
    @{ Cell = "E19"; Value = "-1.06%" },
    @{ Cell = "D20"; Value = "0.1390" },
    @{ Cell = "E20"; Value = "1.39%" },
    @{ Cell = "D22"; Value = "0.04300" },
    @{ Cell = "E22"; Value = "5.14%" },
    @{ Cell = "D23"; Value = "0.001306" },
    @{ Cell = "E23"; Value = "3.17%" },
    @{ Cell = "D24"; Value = "0.004268" },
    @{ Cell = "E24"; Value = "9.07%" },
    @{ Cell = "E25"; Value = "9.61%" },
    @{ Cell = "D26"; Value = "0.0003537" },
    @{ Cell = "E26"; Value = "-5.01%" },
    @{ Cell = "D38"; Value = "0.02753" },
    @{ Cell = "E38"; Value = "14.39%" },
    @{ Cell = "D39"; Value = "0.05619" },
    @{ Cell = "E39"; Value = "7.80%" },
    @{ Cell = "D40"; Value = "0.006297" },
    @{ Cell = "E40"; Value = "-0.14%" },
    @{ Cell = "D41"; Value = "0.007684" },
    @{ Cell = "E41"; Value = "-1.70%" },
    @{ Cell = "D42"; Value = "0.1449" },
    @{ Cell = "E42"; Value = "9.36%" },
    @{ Cell = "D43"; Value = "0.007676" },
    @{ Cell = "E43"; Value = "3.80%" },
    @{ Cell = "D44"; Value = "0.008100" },
    @{ Cell = "E44"; Value = "4.12%" },
    @{ Cell = "D45"; Value = "0.3193" },
    @{ Cell = "E45"; Value = "7.70%" },
    @{ Cell = "E46"; Value = "5.26%" },
    @{ Cell = "E47"; Value = "-0.14%" },
    @{ Cell = "E48"; Value = "35.35%" },
    @{ Cell = "D49"; Value = "0.003998" },
    @{ Cell = "E49"; Value = "-4.90%" },
    @{ Cell = "D50"; Value = "0.00002099" },
    @{ Cell = "E50"; Value = "-0.14%" },
    @{ Cell = "D51"; Value = "0.0001999" },
    @{ Cell = "E51"; Value = "-0.14%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so values like "323.40" / "8.78%" keep their
    # literal, zero-padded string representation instead of being
    # auto-coerced to numeric/percentage types.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Restore the default (unstyled) cell style so no stray formatting
    # is introduced by the temporary text-number-format switch.
    $cell.Style = "Normal"
}
